$d = $word.ActiveDocument

$d.Content.Find.Execute("77-22=", $true, $false, $false, $false, $false, $true, 1, $false, "87-74=", 2) | Out-Null
$d.Content.Find.Execute("2+33=", $true, $false, $false, $false, $false, $true, 1, $false, "95-72=", 2) | Out-Null
$d.Content.Find.Execute("83-32=", $true, $false, $false, $false, $false, $true, 1, $false, "8+22=", 2) | Out-Null
$d.Content.Find.Execute("80-9=", $true, $false, $false, $false, $false, $true, 1, $false, "11-4=", 2) | Out-Null
$d.Content.Find.Execute("90-15=", $true, $false, $false, $false, $false, $true, 1, $false, "6+44=", 2) | Out-Null
$d.Content.Find.Execute("5+29=", $true, $false, $false, $false, $false, $true, 1, $false, "98-0=", 2) | Out-Null
$d.Content.Find.Execute("33-30=", $true, $false, $false, $false, $false, $true, 1, $false, "60-23=", 2) | Out-Null
$d.Content.Find.Execute("18+53=", $true, $false, $false, $false, $false, $true, 1, $false, "96-78=", 2) | Out-Null
$d.Content.Find.Execute("97-8=", $true, $false, $false, $false, $false, $true, 1, $false, "87-32=", 2) | Out-Null
$d.Content.Find.Execute("0+42=", $true, $false, $false, $false, $false, $true, 1, $false, "81+1=", 2) | Out-Null
$d.Content.Find.Execute("92-21=", $true, $false, $false, $false, $false, $true, 1, $false, "6+20=", 2) | Out-Null
$d.Content.Find.Execute("88-81=", $true, $false, $false, $false, $false, $true, 1, $false, "69-27=", 2) | Out-Null
$d.Content.Find.Execute("72-16=", $true, $false, $false, $false, $false, $true, 1, $false, "77-11=", 2) | Out-Null
$d.Content.Find.Execute("1+34=", $true, $false, $false, $false, $false, $true, 1, $false, "65-32=", 2) | Out-Null
$d.Content.Find.Execute("83-36=", $true, $false, $false, $false, $false, $true, 1, $false, "99-77=", 2) | Out-Null
$d.Content.Find.Execute("3+39=", $true, $false, $false, $false, $false, $true, 1, $false, "10-0=", 2) | Out-Null
$d.Content.Find.Execute("22+2=", $true, $false, $false, $false, $false, $true, 1, $false, "85+9=", 2) | Out-Null
$d.Content.Find.Execute("71-32=", $true, $false, $false, $false, $false, $true, 1, $false, "38+48=", 2) | Out-Null
$d.Content.Find.Execute("20+48=", $true, $false, $false, $false, $false, $true, 1, $false, "21-6=", 2) | Out-Null
$d.Content.Find.Execute("78-65=", $true, $false, $false, $false, $false, $true, 1, $false, "15-0=", 2) | Out-Null
$d.Content.Find.Execute("6+38=", $true, $false, $false, $false, $false, $true, 1, $false, "86-84=", 2) | Out-Null
$d.Content.Find.Execute("66-49=", $true, $false, $false, $false, $false, $true, 1, $false, "85-73=", 2) | Out-Null
$d.Content.Find.Execute("72+23=", $true, $false, $false, $false, $false, $true, 1, $false, "0+8=", 2) | Out-Null
$d.Content.Find.Execute("32+41=", $true, $false, $false, $false, $false, $true, 1, $false, "48-27=", 2) | Out-Null
$d.Content.Find.Execute("81-35=", $true, $false, $false, $false, $false, $true, 1, $false, "30-16=", 2) | Out-Null
$d.Content.Find.Execute("11+72=", $true, $false, $false, $false, $false, $true, 1, $false, "56-5=", 2) | Out-Null
$d.Content.Find.Execute("23-12=", $true, $false, $false, $false, $false, $true, 1, $false, "43-17=", 2) | Out-Null
$d.Content.Find.Execute("25+51=", $true, $false, $false, $false, $false, $true, 1, $false, "60-23=", 2) | Out-Null
$d.Content.Find.Execute("87-5=", $true, $false, $false, $false, $false, $true, 1, $false, "58-16=", 2) | Out-Null
$d.Content.Find.Execute("76-26=", $true, $false, $false, $false, $false, $true, 1, $false, "84-68=", 2) | Out-Null
$d.Content.Find.Execute("47-31=", $true, $false, $false, $false, $false, $true, 1, $false, "75-60=", 2) | Out-Null
$d.Content.Find.Execute("74-57=", $true, $false, $false, $false, $false, $true, 1, $false, "51+43=", 2) | Out-Null
$d.Content.Find.Execute("27-10=", $true, $false, $false, $false, $false, $true, 1, $false, "65-32=", 2) | Out-Null
$d.Content.Find.Execute("62-40=", $true, $false, $false, $false, $false, $true, 1, $false, "68-67=", 2) | Out-Null
$d.Content.Find.Execute("45+33=", $true, $false, $false, $false, $false, $true, 1, $false, "79-45=", 2) | Out-Null
$d.Content.Find.Execute("78-7=", $true, $false, $false, $false, $false, $true, 1, $false, "68-48=", 2) | Out-Null
$d.Content.Find.Execute("25+39=", $true, $false, $false, $false, $false, $true, 1, $false, "38+17=", 2) | Out-Null
$d.Content.Find.Execute("2-1=", $true, $false, $false, $false, $false, $true, 1, $false, "44-34=", 2) | Out-Null
$d.Content.Find.Execute("9+82=", $true, $false, $false, $false, $false, $true, 1, $false, "96-66=", 2) | Out-Null
$d.Content.Find.Execute("89-76=", $true, $false, $false, $false, $false, $true, 1, $false, "10+39=", 2) | Out-Null
$d.Content.Find.Execute("59+38=", $true, $false, $false, $false, $false, $true, 1, $false, "96-77=", 2) | Out-Null
$d.Content.Find.Execute("7+1=", $true, $false, $false, $false, $false, $true, 1, $false, "93-50=", 2) | Out-Null
$d.Content.Find.Execute("50-15=", $true, $false, $false, $false, $false, $true, 1, $false, "45-14=", 2) | Out-Null
$d.Content.Find.Execute("8+42=", $true, $false, $false, $false, $false, $true, 1, $false, "77-2=", 2) | Out-Null
$d.Content.Find.Execute("52-26=", $true, $false, $false, $false, $false, $true, 1, $false, "82+1=", 2) | Out-Null
$d.Content.Find.Execute("24+3=", $true, $false, $false, $false, $false, $true, 1, $false, "28+30=", 2) | Out-Null
$d.Content.Find.Execute("52-5=", $true, $false, $false, $false, $false, $true, 1, $false, "78-3=", 2) | Out-Null
$d.Content.Find.Execute("52-7=", $true, $false, $false, $false, $false, $true, 1, $false, "60-52=", 2) | Out-Null
$d.Content.Find.Execute("93-0=", $true, $false, $false, $false, $false, $true, 1, $false, "25+13=", 2) | Out-Null
$d.Content.Find.Execute("0+15=", $true, $false, $false, $false, $false, $true, 1, $false, "65+32=", 2) | Out-Null
$d.Content.Find.Execute("50-24=", $true, $false, $false, $false, $false, $true, 1, $false, "51-42=", 2) | Out-Null
$d.Content.Find.Execute("58-0=", $true, $false, $false, $false, $false, $true, 1, $false, "27+20=", 2) | Out-Null
$d.Content.Find.Execute("48+34=", $true, $false, $false, $false, $false, $true, 1, $false, "63-37=", 2) | Out-Null
$d.Content.Find.Execute("90-29=", $true, $false, $false, $false, $false, $true, 1, $false, "4+3=", 2) | Out-Null
$d.Content.Find.Execute("2+7=", $true, $false, $false, $false, $false, $true, 1, $false, "54+2=", 2) | Out-Null
$d.Content.Find.Execute("57+14=", $true, $false, $false, $false, $false, $true, 1, $false, "61+35=", 2) | Out-Null
$d.Content.Find.Execute("74-23=", $true, $false, $false, $false, $false, $true, 1, $false, "59-18=", 2) | Out-Null
$d.Content.Find.Execute("63+3=", $true, $false, $false, $false, $false, $true, 1, $false, "25+31=", 2) | Out-Null
$d.Content.Find.Execute("39+0=", $true, $false, $false, $false, $false, $true, 1, $false, "10+88=", 2) | Out-Null
$d.Content.Find.Execute("22+73=", $true, $false, $false, $false, $false, $true, 1, $false, "70+9=", 2) | Out-Null
$d.Content.Find.Execute("83-78=", $true, $false, $false, $false, $false, $true, 1, $false, "39+54=", 2) | Out-Null
$d.Content.Find.Execute("1+91=", $true, $false, $false, $false, $false, $true, 1, $false, "66-22=", 2) | Out-Null
$d.Content.Find.Execute("32+11=", $true, $false, $false, $false, $false, $true, 1, $false, "10+45=", 2) | Out-Null
$d.Content.Find.Execute("8+11=", $true, $false, $false, $false, $false, $true, 1, $false, "87-81=", 2) | Out-Null
$d.Content.Find.Execute("3+13=", $true, $false, $false, $false, $false, $true, 1, $false, "64-17=", 2) | Out-Null
$d.Content.Find.Execute("78+16=", $true, $false, $false, $false, $false, $true, 1, $false, "23+12=", 2) | Out-Null
$d.Content.Find.Execute("26-4=", $true, $false, $false, $false, $false, $true, 1, $false, "69-40=", 2) | Out-Null
$d.Content.Find.Execute("70+24=", $true, $false, $false, $false, $false, $true, 1, $false, "51+38=", 2) | Out-Null
$d.Content.Find.Execute("99-5=", $true, $false, $false, $false, $false, $true, 1, $false, "6+5=", 2) | Out-Null
$d.Content.Find.Execute("27+64=", $true, $false, $false, $false, $false, $true, 1, $false, "6+77=", 2) | Out-Null
$d.Content.Find.Execute("49+0=", $true, $false, $false, $false, $false, $true, 1, $false, "53+21=", 2) | Out-Null
$d.Content.Find.Execute("64-30=", $true, $false, $false, $false, $false, $true, 1, $false, "6+28=", 2) | Out-Null
$d.Content.Find.Execute("21-13=", $true, $false, $false, $false, $false, $true, 1, $false, "6+62=", 2) | Out-Null
$d.Content.Find.Execute("45+12=", $true, $false, $false, $false, $false, $true, 1, $false, "61+33=", 2) | Out-Null
$d.Content.Find.Execute("19+25=", $true, $false, $false, $false, $false, $true, 1, $false, "82-35=", 2) | Out-Null
$d.Content.Find.Execute("72+22=", $true, $false, $false, $false, $false, $true, 1, $false, "67-8=", 2) | Out-Null
$d.Content.Find.Execute("1+65=", $true, $false, $false, $false, $false, $true, 1, $false, "73-4=", 2) | Out-Null
$d.Content.Find.Execute("96-58=", $true, $false, $false, $false, $false, $true, 1, $false, "50-22=", 2) | Out-Null
$d.Content.Find.Execute("35+57=", $true, $false, $false, $false, $false, $true, 1, $false, "6+75=", 2) | Out-Null
$d.Content.Find.Execute("79-1=", $true, $false, $false, $false, $false, $true, 1, $false, "25-15=", 2) | Out-Null
$d.Content.Find.Execute("6+76=", $true, $false, $false, $false, $false, $true, 1, $false, "45-44=", 2) | Out-Null
$d.Content.Find.Execute("91-13=", $true, $false, $false, $false, $false, $true, 1, $false, "54+0=", 2) | Out-Null
$d.Content.Find.Execute("23-3=", $true, $false, $false, $false, $false, $true, 1, $false, "28+31=", 2) | Out-Null
$d.Content.Find.Execute("82-30=", $true, $false, $false, $false, $false, $true, 1, $false, "60-10=", 2) | Out-Null
$d.Content.Find.Execute("79-52=", $true, $false, $false, $false, $false, $true, 1, $false, "6+7=", 2) | Out-Null
$d.Content.Find.Execute("53-27=", $true, $false, $false, $false, $false, $true, 1, $false, "89-11=", 2) | Out-Null
$d.Content.Find.Execute("43-21=", $true, $false, $false, $false, $false, $true, 1, $false, "72-55=", 2) | Out-Null
$d.Content.Find.Execute("64-19=", $true, $false, $false, $false, $false, $true, 1, $false, "51-20=", 2) | Out-Null
$d.Content.Find.Execute("16+41=", $true, $false, $false, $false, $false, $true, 1, $false, "64-21=", 2) | Out-Null
$d.Content.Find.Execute("77-71=", $true, $false, $false, $false, $false, $true, 1, $false, "94-75=", 2) | Out-Null
$d.Content.Find.Execute("59-37=", $true, $false, $false, $false, $false, $true, 1, $false, "27+33=", 2) | Out-Null
$d.Content.Find.Execute("35+45=", $true, $false, $false, $false, $false, $true, 1, $false, "43-37=", 2) | Out-Null
$d.Content.Find.Execute("66-15=", $true, $false, $false, $false, $false, $true, 1, $false, "69-17=", 2) | Out-Null
$d.Content.Find.Execute("2+30=", $true, $false, $false, $false, $false, $true, 1, $false, "34-7=", 2) | Out-Null
$d.Content.Find.Execute("55+27=", $true, $false, $false, $false, $false, $true, 1, $false, "2+81=", 2) | Out-Null
$d.Content.Find.Execute("31+38=", $true, $false, $false, $false, $false, $true, 1, $false, "65-3=", 2) | Out-Null
$d.Content.Find.Execute("77-20=", $true, $false, $false, $false, $false, $true, 1, $false, "37+27=", 2) | Out-Null
$d.Content.Find.Execute("52+3=", $true, $false, $false, $false, $false, $true, 1, $false, "41+32=", 2) | Out-Null
$d.Content.Find.Execute("23+41=", $true, $false, $false, $false, $false, $true, 1, $false, "7-2=", 2) | Out-Null
$d.Content.Find.Execute("43+26=", $true, $false, $false, $false, $false, $true, 1, $false, "74+0=", 2) | Out-Null
